# "Antworten mischen" (mix_answers) column added for SC questions.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column BS: header + two data rows -------------------------------
# Header cell BS1 gets the same "header" look as the other header cells
# (bold font on a yellow fill) but only with a left/right border instead of
# the full box border used by the rest of row 1.
$hdr = $ws.Range("BS1")
$hdr.Value = "mix_answers"
$hdr.Font.Bold = $true
$hdr.Interior.Color = 65535

$hdr.Borders.Item(7).LineStyle = 1    # xlEdgeLeft  / xlContinuous
$hdr.Borders.Item(10).LineStyle = 1   # xlEdgeRight / xlContinuous

# Data rows: 0 = don't mix answers, 1 = mix answers
$ws.Range("BS2").Value = 0
$ws.Range("BS3").Value = 1

# --- View state: scroll right and move the selection ----------------------
$win = $excel.ActiveWindow
$win.ScrollColumn = 45   # AS
$win.ScrollRow = 1

$ws.Range("BR9").Select()
